# uk_local_authorities.xlsx - "Added more details on combined authorities
# And matched to constituent councils"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Column F ("overlapping-la"): tag the constituent councils of each
#    combined authority with the CA's BS-6879 abbreviation.
# ---------------------------------------------------------------------------
$overlap = @{
    26  = "SCR";
    27  = "GMCA";
    34  = "WYCA";
    44  = "GMCA";
    63  = "WYCA";
    77  = "TVCA";
    84  = "SCR";
    89  = "NECA";
    117 = "NECA";
    131 = "LCR";
    147 = "TVCA";
    162 = "WYCA";
    164 = "LCR";
    170 = "WYCA";
    177 = "LCR";
    182 = "GMCA";
    185 = "TVCA";
    195 = "NTCA";
    202 = "NTCA";
    215 = "NTCA";
    222 = "GMCA";
    232 = "TVCA";
    233 = "GMCA";
    244 = "SCR";
    263 = "LCR";
    267 = "SCR";
    268 = "LCR";
    272 = "GMCA";
    274 = "GMCA";
    276 = "NECA";
    296 = "TVCA";
    298 = "NECA";
    305 = "GMCA";
    320 = "GMCA";
    340 = "GMCA";
    343 = "WYCA";
    357 = "LCR"
}
foreach ($row in $overlap.Keys) {
    $ws.Cells.Item($row, 6).Value = $overlap[$row]
}

# ---------------------------------------------------------------------------
# 2) Northern Ireland / Scotland / Wales blocks (rows 370-459): add
#    local-authority-type (D) + local-authority-type-name (E), and fill in
#    the missing region (C) for the rows that didn't have it yet.
# ---------------------------------------------------------------------------
for ($row = 370; $row -le 405; $row++) {
    if ($row -ge 380 -and $row -le 404) {
        $ws.Cells.Item($row, 3).Value = "Northern Ireland"
    }
    $ws.Cells.Item($row, 4).Value = "NID"
    $ws.Cells.Item($row, 5).Value = "NI District Council"
}

for ($row = 406; $row -le 437; $row++) {
    $ws.Cells.Item($row, 4).Value = "SCO"
    $ws.Cells.Item($row, 5).Value = "Scottish Unitary Council"
}

for ($row = 438; $row -le 459; $row++) {
    $ws.Cells.Item($row, 4).Value = "WPA"
    $ws.Cells.Item($row, 5).Value = "Welsh Unitary Council"
}

# ---------------------------------------------------------------------------
# 3) GSS codes for the city-region combined authorities (gss-code, col O)
# ---------------------------------------------------------------------------
$ws.Cells.Item(169, 15).Value = "E47000004"   # Liverpool City Region
$ws.Cells.Item(256, 15).Value = "E47000002"   # Sheffield City Region

# ---------------------------------------------------------------------------
# 4) BS-6879 codes (col V) for combined authorities - correct the rows that
#    had been matched against the wrong/legacy code, and add the code for
#    the new NECA / WYCA / NTCA rows at the bottom of the sheet.
# ---------------------------------------------------------------------------
$ws.Cells.Item(349, 22).Value = "NECA"
$ws.Cells.Item(350, 22).Value = "WYCA"
$ws.Cells.Item(351, 22).Value = "NTCA"

$ws.Cells.Item(461, 16).Value = "E47000005"   # archaic-gss-code for NECA
$ws.Cells.Item(461, 22).Value = "NECA"
$ws.Cells.Item(462, 22).Value = "WYCA"
$ws.Cells.Item(463, 22).Value = "NTCA"

# ---------------------------------------------------------------------------
# 5) Cosmetic sheet-view tidy up: scroll back to the top-left, re-freeze the
#    header row at A2, and widen column E to fit its new content.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 29.5
$ws.Range("E10").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 6) Add the blank "Sheet1" tab after the data sheet, then re-activate the
#    original sheet so it stays the one shown when the file is opened.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Range("C17").Select()
$ws.Activate()
